$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update column C (Förändrad) for rows 2-5 from 45174 to 45175
$ws.Range("C2:C5").Value = 45175
